$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("Cell Topcon 183mm")
$ws.Range("B12").NumberFormat = "@"
$ws.Range("B12").Value = "0.298"

$ws = $wb.Worksheets.Item("Silver Rear_side")
$ws.Range("B12").NumberFormat = "@"
$ws.Range("B12").Value = "5,335"

$ws = $wb.Worksheets.Item("Silver Busbar front-side")
$ws.Range("B12").NumberFormat = "@"
$ws.Range("B12").Value = "7,987"

$ws = $wb.Worksheets.Item("Silver finger front-side")
$ws.Range("B12").NumberFormat = "@"
$ws.Range("B12").Value = "8,037"

$ws = $wb.Worksheets.Item("USD_CNY")
$ws.Range("B12").NumberFormat = "@"
$ws.Range("B12").Value = "7.2567"
